$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-CellText $ws "D2" "68.716.95"
Set-CellText $ws "E2" "  +2.34%  "
Set-CellText $ws "D3" "2.528.42"
Set-CellText $ws "E4" "  +0.03%  "
Set-CellText $ws "D5" "594.43"
Set-CellText $ws "E5" "  +2.03%  "
Set-CellText $ws "D6" "177.44"
Set-CellText $ws "E6" "  +1.54%  "
Set-CellText $ws "E7" "  +0.02%  "
Set-CellText $ws "D8" "0.521"
Set-CellText $ws "E8" "  +1.69%  "
Set-CellText $ws "D9" "2.528.18"
Set-CellText $ws "D10" "0.147"
Set-CellText $ws "E10" "  +6.37%  "
Set-CellText $ws "E11" "  -1.03%  "
Set-CellText $ws "E12" "  +1.08%  "
Set-CellText $ws "D13" "0.340"
Set-CellText $ws "E13" "  +1.76%  "
Set-CellText $ws "D14" "2.990.90"
Set-CellText $ws "E14" "  +2.74%  "
Set-CellText $ws "D15" "26.25"
Set-CellText $ws "E15" "  +3.42%  "
Set-CellText $ws "D16" "68.262.14"
Set-CellText $ws "E16" "  +1.90%  "
Set-CellText $ws "D17" "0.0000171"
Set-CellText $ws "E17" "  +0.92%  "
Set-CellText $ws "D18" "2.519.95"
Set-CellText $ws "E18" "  +1.91%  "
Set-CellText $ws "D19" "11.12"
Set-CellText $ws "E19" "  +1.85%  "
Set-CellText $ws "D20" "7.52"
Set-CellText $ws "E20" "  +0.77%  "
Set-CellText $ws "D21" "353.03"
Set-CellText $ws "E21" "  +1.42%  "
Set-CellText $ws "E22" "  +4.72%  "
Set-CellText $ws "D23" "1.00"
Set-CellText $ws "E23" "  +0.05%  "
Set-CellText $ws "D24" "70.99"
Set-CellText $ws "E24" "  +2.34%  "
Set-CellText $ws "D25" "4.26"
Set-CellText $ws "E25" "  +1.58%  "
Set-CellText $ws "E26" "  -4.94%  "
Set-CellText $ws "D27" "9.00"
Set-CellText $ws "E27" "  -2.27%  "
Set-CellText $ws "D28" "2.690.84"
Set-CellText $ws "E28" "  +3.70%  "
Set-CellText $ws "D29" "0.995"
Set-CellText $ws "E29" "  -0.49%  "
Set-CellText $ws "B30" "PEPE"
Set-CellText $ws "C30" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-CellText $ws "D30" "0.0₃0895"
Set-CellText $ws "E30" "  -0.49%  "
Set-CellText $ws "B31" "Bittensor"
Set-CellText $ws "C31" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-CellText $ws "D31" "509.08"
Set-CellText $ws "E31" "  +1.96%  "
Set-CellText $ws "D32" "7.81"
Set-CellText $ws "E32" "  +1.08%  "
Set-CellText $ws "E33" "  +1.84%  "
Set-CellText $ws "E34" "  +1.28%  "
Set-CellText $ws "D35" "1.00"
Set-CellText $ws "E35" "  +0.05%  "
Set-CellText $ws "B36" "Kaspa"
Set-CellText $ws "C36" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-CellText $ws "D36" "0.120"
Set-CellText $ws "E36" "  +0.01%  "
Set-CellText $ws "B37" "Monero"
Set-CellText $ws "C37" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-CellText $ws "D37" "162.94"
Set-CellText $ws "E37" "  +0.90%  "
Set-CellText $ws "B38" "EthereumClassic"
Set-CellText $ws "C38" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-CellText $ws "D38" "18.43"
Set-CellText $ws "E38" "  +1.47%  "
Set-CellText $ws "B39" "WhiteBITCoin"
Set-CellText $ws "C39" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-CellText $ws "D39" "18.68"
Set-CellText $ws "E39" "  +0.02%  "
Set-CellText $ws "E40" "  +5.28%  "
Set-CellText $ws "E41" "  -0.48%  "
Set-CellText $ws "E42" "  +0.01%  "
Set-CellText $ws "D43" "4.85"
Set-CellText $ws "E43" "  +0.62%  "
Set-CellText $ws "E44" "  -0.18%  "
Set-CellText $ws "D45" "2.42"
Set-CellText $ws "E45" "  +1.38%  "
Set-CellText $ws "D46" "153.30"
Set-CellText $ws "E46" "  +7.47%  "
Set-CellText $ws "E47" "  +2.77%  "
Set-CellText $ws "D48" "0.522"
Set-CellText $ws "E48" "  +2.44%  "
Set-CellText $ws "D49" "0.0₆0260"
Set-CellText $ws "E49" "  +1.62%  "
Set-CellText $ws "E50" "  +2.67%  "
